$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the old extra rows (7-10) - the validated/cleaned data set only
#    keeps 6 data rows now ("se valida que si estan vacios no haga nada").
# ---------------------------------------------------------------------------
$ws.Range("A7:E10").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Replace the 6 remaining rows with the new, validated pedimentos data.
#    Columns A and C hold numeric/date-looking text ("135330613002049",
#    "12/19/2013", ...) so they are forced to Text (NumberFormat "@") before
#    the value is written, otherwise Excel's smart entry would turn them
#    into numbers/dates. Columns B, D and E are unambiguous text already.
# ---------------------------------------------------------------------------
$data = @(
    @("135330613002049", "AICM", "12/19/2013", "13-802240", "CANCUN"),
    @("124735222003487", "AICM", "12/24/2012", "12-03350",  "AICM"),
    @("124735222003488", "AICM", "12/24/2012", "12-03288",  "AICM"),
    @("134735223000062", "AICM", "02/15/2013", "13-00085",  "AICM"),
    @("134735223004059", "AICM", "12/26/2013", "13-03918",  "AICM"),
    @("135330613000523", "AICM", "03/27/2013", "13-800505", "CANCUN")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 1
    $row = $data[$i]

    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = $row[2]

    $ws.Cells.Item($r, 4).Value = $row[3]

    $ws.Cells.Item($r, 5).Value = $row[4]
}

# ---------------------------------------------------------------------------
# 3. Formatting per column, matching the cleaned-up template:
#    - Column A: hidden pedimento numbers (";;" format, Arial 10)
#    - Column C: dates, left aligned; row 1 is a text-style header
#    - D1/E1/E6 reuse the hidden Arial 10 style as well
# ---------------------------------------------------------------------------
for ($r = 1; $r -le 6; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = ";;"
    $cellA.Font.Name = "Arial"
    $cellA.Font.Size = 10
    $cellA.Font.Bold = $false
}

# Row 1 (header-like row): C1 is a left-aligned text cell, Arial 10
$c1 = $ws.Cells.Item(1, 3)
$c1.NumberFormat = "@"
$c1.Font.Name = "Arial"
$c1.Font.Size = 10
$c1.HorizontalAlignment = -4131

# Rows 2-5: dates, left aligned, default font
foreach ($r in 2..5) {
    $c = $ws.Cells.Item($r, 3)
    $c.NumberFormat = "m/d/yyyy"
    $c.HorizontalAlignment = -4131
}

# Row 6: date, default alignment
$c6 = $ws.Cells.Item(6, 3)
$c6.NumberFormat = "m/d/yyyy"

# D1 / E1 / E6 reuse the hidden Arial 10 style
foreach ($ref in @("D1", "E1", "E6")) {
    $c = $ws.Range($ref)
    $c.NumberFormat = ";;"
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.Font.Bold = $false
}

# ---------------------------------------------------------------------------
# 4. Column widths
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.28515625
$ws.Columns.Item(2).ColumnWidth = 41.140625
$ws.Columns.Item(3).ColumnWidth = 24.7109375
$ws.Columns.Item(4).ColumnWidth = 15.5703125
$ws.Columns.Item(5).ColumnWidth = 12.85546875

# ---------------------------------------------------------------------------
# 5. Selection moves to A2, page orientation is set to portrait
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$ws.PageSetup.Orientation = 1

Write-Output "done"
